$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 143; this shifts rows 143:167 down to 144:168
# and extends the sheet dimension to A1:R168 automatically, carrying over
# row formatting (e.g. the date style on column D).
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with the new weekly data record.
$ws.Range("A143").Value = 9
$ws.Range("B143").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C143").Value = "Metropolitana"
$ws.Range("D143").Value = 44505
$ws.Range("E143").Value = 13
$ws.Range("F143").Value = 100112026
$ws.Range("G143").Value = "Haba"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 52
$ws.Range("K143").Value = 7000
$ws.Range("L143").Value = 8000
$ws.Range("M143").Value = 7500
$ws.Range("N143").Value = "$/saco 25 kilos"
$ws.Range("O143").Value = "Región Metropolitana"
$ws.Range("P143").Value = 300
$ws.Range("Q143").Value = 25
$ws.Range("R143").Value = "Hortaliza"
